# Refresh of the Query1 / ExternalData_1 queryTable against its source:
# one row ("ALLEGRETTO (B7981027)") no longer present, and several
# Progress values changed as a result of the refreshed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ALLEGRETTO (B7981027)" trial row (row 9) dropped out of the refreshed
# query result entirely - remove it, shifting all following rows up.
$ws.Range("A9").EntireRow.Delete()

# Updated Progress values coming from the refreshed query result.
$ws.Range("B8").Value = 12    # HORIZON OLE
$ws.Range("B9").Value = 25    # ALLEGRETTO-LTE (B7981028)
$ws.Range("B15").Value = 0    # ALPINE
$ws.Range("B16").Value = 0    # REMASTER (CLOU)

# Keep the hidden ExternalData_1 defined name (which tracks the query's
# result range) in sync with the now-smaller table extent.
$wb.Names("ExternalData_1").RefersTo = "=Sheet1!`$A`$1:`$B`$16"
